$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rotate C1/D1/E1 -> C1=prediction, D1=rejection-f, E1=max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes text (previously D2's value), D2 stays the same, E2 becomes numeric 1
$ws.Range("C2").Value = "c__Fusobacteriia"
$ws.Range("E2").Value = 1

# Row 3: C3 becomes text (previously D3's value), D3 stays the same, E3 becomes numeric 1
$ws.Range("C3").Value = "c__Fusobacteriia"
$ws.Range("E3").Value = 1
